$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared strings must be interned in the same order as in the target
# workbook: "Jon" first, then "Eliminar con boton...", then "Crear con
# forularios...". Write the B column (names) first, then fill in C5 before
# C4 so the shared-string table ends up in that exact order.
$ws.Range("B4").Value = "Jon"
$ws.Range("B5").Value = "Gaizka"
$ws.Range("C5").Value = "Eliminar con botón desde detalles: Proyecto, Tarea, Empleado"
$ws.Range("C4").Value = "Crear con forularios: Proyecto, Tarea, Empleado"

# Apply the date number format *before* assigning the value, otherwise
# Excel auto-applies its own default date format first (creating an extra,
# unused style entry) and only then gets overridden by our explicit format.
$ws.Range("D4").NumberFormat = "d-mmm"
$ws.Range("D4").Value = (Get-Date -Year 2025 -Month 4 -Day 19 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("D5").NumberFormat = "d-mmm"
$ws.Range("D5").Value = (Get-Date -Year 2025 -Month 4 -Day 20 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)

$ws.Range("D5").Select()
